$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Increase DMG for Slugs (row 38, ammo_12x76_zhekan - Perf type)
$ws.Range("H38").Value = 2.7

# Increase DMG for Buckshot (row 39, ammo_12x70_buck - DMG type)
$ws.Range("H39").Formula = "=9*0.42"

# Update the active selection to match the recorded workbook state
$ws.Range("J27").Select()
